# Eggerland.xlsx - "many more levels" commit
# Adds 21 new timing rows (55-75) to Sheet1, three new lookup strings used
# by column A, and a handful of scratch formulas in columns I/J that were
# present in the author's working area while adding the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data rows 55-75 (columns A-D), continuing the existing table ---
# Column A uses the same lookup-text convention as the existing rows;
# column D repeats the existing "C-B" timing-delta formula.

$ws.Range("A55").Value = "Get key"
$ws.Range("B55").Value = 92262
$ws.Range("C55").Value = 90460
$ws.Range("D55").Formula = "=C55-B55"

$ws.Range("A56").Value = "Green Lolo"
$ws.Range("B56").Value = 92927
$ws.Range("C56").Value = 91126
$ws.Range("D56").Formula = "=C56-B56"

$ws.Range("A57").Value = "Get key"
$ws.Range("B57").Value = 93439
$ws.Range("C57").Value = 91638
$ws.Range("D57").Formula = "=C57-B57"

$ws.Range("A58").Value = "Get key"
$ws.Range("B58").Value = 96846
$ws.Range("C58").Value = 95046
$ws.Range("D58").Formula = "=C58-B58"

$ws.Range("A59").Value = "Level begin"
$ws.Range("B59").Value = 97744
$ws.Range("C59").Value = 95821
$ws.Range("D59").Formula = "=C59-B59"

# Row 60 has no label in column A (matches the source data).
$ws.Range("B60").Value = 98315
$ws.Range("C60").Value = 96392
$ws.Range("D60").Formula = "=C60-B60"
$ws.Range("J60").Formula = "=5/6"

$ws.Range("A61").Value = "Get treasure"
$ws.Range("B61").Value = 98587
$ws.Range("C61").Value = 96664
$ws.Range("D61").Formula = "=C61-B61"

$ws.Range("A62").Value = "Level begin"
$ws.Range("B62").Value = 99181
$ws.Range("C62").Value = 97136
$ws.Range("D62").Formula = "=C62-B62"

$ws.Range("A63").Value = "Level begin"
$ws.Range("B63").Value = 99784
$ws.Range("C63").Value = 97739
$ws.Range("D63").Formula = "=C63-B63"
$ws.Range("I63").Formula = "=123*3"

$ws.Range("A64").Value = "Level scroll"
$ws.Range("B64").Value = 102281
$ws.Range("C64").Value = 100238
$ws.Range("D64").Formula = "=C64-B64"
$ws.Range("I64").Formula = "=369+2298"

$ws.Range("A65").Value = "Level scroll"
$ws.Range("B65").Value = 102435
$ws.Range("C65").Value = 100392
$ws.Range("D65").Formula = "=C65-B65"
$ws.Range("I65").Formula = "=I64/60"

$ws.Range("A66").Value = "Get key"
$ws.Range("B66").Value = 104783
$ws.Range("C66").Value = 102732
$ws.Range("D66").Formula = "=C66-B66"

$ws.Range("A67").Value = "Get key"
$ws.Range("B67").Value = 105147
$ws.Range("C67").Value = 103096
$ws.Range("D67").Formula = "=C67-B67"

$ws.Range("A68").Value = "Level appear"
$ws.Range("B68").Value = 105798
$ws.Range("C68").Value = 103624
$ws.Range("D68").Formula = "=C68-B68"

$ws.Range("A69").Value = "Get treasure"
$ws.Range("B69").Value = 106665
$ws.Range("C69").Value = 104491
$ws.Range("D69").Formula = "=C69-B69"

$ws.Range("A70").Value = "Level appear"
$ws.Range("B70").Value = 107257
$ws.Range("C70").Value = 104959
$ws.Range("D70").Formula = "=C70-B70"

$ws.Range("A71").Value = "Get key"
$ws.Range("B71").Value = 110558
$ws.Range("C71").Value = 108260
$ws.Range("D71").Formula = "=C71-B71"

$ws.Range("A72").Value = "Get key"
$ws.Range("B72").Value = 112188
$ws.Range("C72").Value = 109890
$ws.Range("D72").Formula = "=C72-B72"

$ws.Range("A73").Value = "Get key"
$ws.Range("B73").Value = 113158
$ws.Range("C73").Value = 110860
$ws.Range("D73").Formula = "=C73-B73"

$ws.Range("A74").Value = "Get key"
$ws.Range("B74").Value = 115248
$ws.Range("C74").Value = 112951
$ws.Range("D74").Formula = "=C74-B74"

$ws.Range("A75").Value = "Get key"
$ws.Range("B75").Value = 116008
$ws.Range("C75").Value = 113711
$ws.Range("D75").Formula = "=C75-B75"

# --- View state: scroll so row 59 is at the top, select the next blank row ---
$win = $excel.ActiveWindow
$ws.Range("B76").Select() | Out-Null
$win.ScrollRow = 59
$win.ScrollColumn = 1
